# RICEFW_STATUS.xlsx update
# - "ID Card" row (row 6) added, filled in first (drives new shared-string
#   indices for ID Card / SRMU student id card / psotGoLive / WIP in that order)
# - Row4 status Done -> Completed
# - Row5 status InProgress -> WIP (reuses the WIP string just created)
# - Row5 gets its Type/area/module (B5:D5) filled in (copy of row4's)
# - "Cash Book Summary" row (row 7) added
# - Row6's Type cell (B6) corrected to "REports" last
# - Selection/view updated to E5, no more split scroll to H1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: fill in the previously-empty Type/area/module cells ---
$ws.Cells.Item(5, 2).Value = "Reports"             # B5
$ws.Cells.Item(5, 3).Value = "campus"              # C5
$ws.Cells.Item(5, 4).Value = "Student financials"  # D5

# --- Row 6 (new): ID Card ---
$ws.Cells.Item(6, 3).Value  = "campus"                               # C6
$ws.Cells.Item(6, 4).Value  = "Student financials"                   # D6
$ws.Cells.Item(6, 6).Value  = "N"                                    # F6
$ws.Cells.Item(6, 7).Value  = "Y"                                    # G6
$ws.Cells.Item(6, 8).Value  = "GAP"                                  # H6
$ws.Cells.Item(6, 9).Value  = "ID Card"                              # I6
$ws.Cells.Item(6, 10).Value = "SRMU student id card"                 # J6
$ws.Cells.Item(6, 12).Value = "psotGoLive"                           # L6
$ws.Cells.Item(6, 13).Value = "WIP"                                  # M6
$ws.Cells.Item(6, 14).Value = "Tushar"                                # N6
$ws.Cells.Item(6, 15).Value = "NO"                                   # O6

# --- Status updates on existing rows ---
$ws.Cells.Item(4, 13).Value = "Completed"   # M4: Done -> Completed
$ws.Cells.Item(5, 13).Value = "WIP"         # M5: InProgress -> WIP

# --- Row 7 (new): Cash Book Summary ---
$ws.Cells.Item(7, 2).Value  = "Reports"                                                              # B7
$ws.Cells.Item(7, 3).Value  = "campus"                                                                # C7
$ws.Cells.Item(7, 4).Value  = "Student financials"                                                    # D7
$ws.Cells.Item(7, 6).Value  = "N"                                                                     # F7
$ws.Cells.Item(7, 7).Value  = "Y"                                                                     # G7
$ws.Cells.Item(7, 8).Value  = "GAP"                                                                   # H7
$ws.Cells.Item(7, 9).Value  = "Cash Book Summary"                                                     # I7
$ws.Cells.Item(7, 10).Value = "To print cash book between 2 particular dates (i.e. from date and to date)"  # J7
$ws.Cells.Item(7, 12).Value = "post GoLive "                                                          # L7
$ws.Cells.Item(7, 13).Value = "Completed"                                                              # M7
$ws.Cells.Item(7, 14).Value = "Tushar"                                                                 # N7
$ws.Cells.Item(7, 15).Value = "NO"                                                                    # O7

# --- Row 6 Type column fixed last (creates the final new shared string) ---
$ws.Cells.Item(6, 2).Value = "REports"   # B6

# --- View/selection update ---
$ws.Range("E5").Select() | Out-Null
